$wb = $excel.ActiveWorkbook

# --- Login sheet: remove static credentials' mailto hyperlinks, update the ---
# --- password test value, and rest the hyperlink-style font on the cells  ---
$wsLogin = $wb.Worksheets.Item("Login")
$wsLogin.Range("A2").Hyperlinks.Delete()
$wsLogin.Range("B2").Hyperlinks.Delete()
$wsLogin.Range("B2").Value = "test123"
$wsLogin.Range("A2:B2").Style = "Hyperlink"
$wsLogin.Range("B2").Select()

# --- ChangePassword sheet: same cleanup, new static password values ---
$wsChangePassword = $wb.Worksheets.Item("ChangePassword")
$wsChangePassword.Range("A2").Hyperlinks.Delete()
$wsChangePassword.Range("B2").Hyperlinks.Delete()
$wsChangePassword.Range("C2").Hyperlinks.Delete()
$wsChangePassword.Range("A2").Value = "test1234"
$wsChangePassword.Range("B2").Value = "test123"
$wsChangePassword.Range("C2").Value = "test123"
$wsChangePassword.Range("A2:C2").Style = "Hyperlink"
$wsChangePassword.Range("D3").Select()
